# Adds a new "master" worksheet at the end of the workbook summarising the
# Global Corporate Average Cumulative Default Rates as plain numbers
# (rather than the pre-formatted text used on "adjust" / "no_adjust"),
# highlighting the AAA row and restyling the figures used in the result
# section (commit: "Added additional figures for result section").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet by copying "no_adjust" to the end of the tab strip.
#    Copying (rather than Worksheets.Add()) carries over the column width,
#    row heights and existing cell styles (A1/A3 headers, the 1..15 rating
#    scale in row 3, etc.) so we only need to touch what actually changes.
# ---------------------------------------------------------------------------
$source = $wb.Worksheets.Item("no_adjust")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "master"

# ---------------------------------------------------------------------------
# 2. Wipe the inherited (text) figures in B4:P13 - the master sheet gets its
#    own numeric figures below.
# ---------------------------------------------------------------------------
$ws3.Range("B4:P13").ClearContents()

# ---------------------------------------------------------------------------
# 3. Write the new numeric figures (rows 4-9, cols B-K).
# ---------------------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K")

$row4 = @(0.00001, 0.03, 0.099, 0.189, 0.27, 0.369, 0.40599999999999997, 0.45629629629629631, 0.52499999999999991, 0.58099999999999996)
$row5 = @(0.02, 0.06, 0.11, 0.21, 0.3, 0.41, 0.49, 0.56000000000000005, 0.63, 0.7)
$row6 = @(0.05, 0.13, 0.22, 0.33, 0.46, 0.6, 0.76, 0.9, 1.05, 1.2)
$row7 = @(0.16, 0.43, 0.75, 1.1399999999999999, 1.54, 1.94, 2.27, 2.61, 2.94, 3.24)
$row8 = @(0.63, 1.93, 3.46, 4.99, 6.43, 7.75, 8.89, 9.9, 10.82, 11.64)
$row9 = @(3.34, 7.8, 11.75, 14.89, 17.350000000000001, 19.36, 20.99, 22.31, 23.5, 24.62)

$data = @{ 4 = $row4; 5 = $row5; 6 = $row6; 7 = $row7; 8 = $row8; 9 = $row9 }

foreach ($r in @(4,5,6,7,8,9)) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws3.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 4. Formatting.
#    - Row 4 (AAA) is called out in red, sized up to 12pt, centred, "0.00".
#    - Rows 5-9 (AA..B) are 12pt, centred; plain General format except for
#      F5, which keeps the "0.00" look.
# ---------------------------------------------------------------------------
$ws3.Range("B4:K4").Font.Size = 12
$ws3.Range("B4:K4").Font.Color = 255
$ws3.Range("B4:K4").NumberFormat = "0.00"
$ws3.Range("B4:K4").HorizontalAlignment = -4108

$ws3.Range("B5:K9").Font.Size = 12
$ws3.Range("B5:K9").HorizontalAlignment = -4108

$ws3.Range("F5").NumberFormat = "0.00"

# Blank-but-styled cells: the unused right-hand columns (L:P) for rows 4-9
# and the whole of rows 10-13 keep the workbook's existing "0.00" / left
# aligned look (same style already used elsewhere in the workbook).
$ws3.Range("L4:P9").NumberFormat = "0.00"
$ws3.Range("L4:P9").HorizontalAlignment = -4131

$ws3.Range("B10:P13").NumberFormat = "0.00"
$ws3.Range("B10:P13").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 5. Row heights - rows 6-9 grow slightly to fit the bigger 12pt font.
# ---------------------------------------------------------------------------
foreach ($r in @(6,7,8,9)) {
    $ws3.Rows.Item($r).RowHeight = 15.5
}

# ---------------------------------------------------------------------------
# 6. Selection / active tab - "master" becomes the active sheet, with M8
#    selected; "adjust" (previously active) loses tabSelected automatically.
# ---------------------------------------------------------------------------
$ws3.Range("M8").Select()
